$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove computed/value content for row 246 (28.10.2020) - clear formulas/values
$ws.Range("B246").ClearContents()
$ws.Range("D246").ClearContents()
$ws.Range("H246").ClearContents()
$ws.Range("J246").ClearContents()
$ws.Range("K246").ClearContents()

# Update frozen-pane scroll position and active selection to match new view state
$ws.Range("F247").Select()
$excel.ActiveWindow.ScrollRow = 222
